# Ind_Int_EEUU.xlsx — update historical "Produccion industrial (Indice total)"
# values in column D (revised data for an updated source upload).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("C_7")

$updates = @{
    23 = 102.9141
    25 = 101.66070000000001
    26 = 102.4696
    27 = 103.0849
    28 = 103.4627
    29 = 104.8612
    30 = 103.34439999999999
    31 = 103.92619999999999
    32 = 101.9909
    33 = 100.8862
    34 = 103.2771
    35 = 102.0634
    36 = 101.5271
    37 = 100.6437
    38 = 100.6087
    39 = 100.1228
    40 = 98.75
    41 = 101.89190000000001
    42 = 100.1131
    43 = 100.80200000000001
    44 = 98.150300000000001
    45 = 96.961500000000001
    46 = 98.701999999999998
    47 = 95.888599999999997
    48 = 98.966499999999996
    49 = 97.869399999999999
    50 = 96.284700000000001
    51 = 96.455699999999993
    52 = 95.703000000000003
    53 = 98.2239
    54 = 95.093900000000005
    55 = 92.993700000000004
    56 = 84.979799999999997
    57 = 82.676699999999997
    58 = 97.961399999999998
    59 = 101.9971
    60 = 101.2492
}

foreach ($row in $updates.Keys) {
    $ws.Range("D$row").Value = $updates[$row]
}
